$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slurry")

# Update the delta pH formula/value for the 3.4 kg/t acid dose rows (D4 and D8)
$ws.Range("D4").Formula = "=7.9-0.8187"
$ws.Range("D8").Formula = "=7.9-0.8187"

# Update the active cell selection on the sheet, as recorded in the saved view
$ws.Activate()
$ws.Range("E8").Select()
